# Generate Report for Handoff
# Adds two new "dependency" file rows (png pair) to each of the three
# sheets (Overview / zh-cn / de-de), refreshes the hand-off artefact
# names + timestamps for the existing ".md"/".xlf" row, and wires up
# hyperlinks for every file-name cell, matching the authoring commit
# "Generate Report for Handoff".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Shared literals
# ---------------------------------------------------------------------
$oldMdUuid   = "92ee2ffd-a547-4a8f-a31c-0bd2f950989f"
$newMdUuid   = "389c4df5-9bb7-413b-b08f-2776d8518a66"
$oldXlfHash  = "465aa127ac377bb7970b35a85db41248b6fbed9e"
$newXlfHash  = "f21aa4a2eda604d85fc8e8d7d8c3c138fa2649e2"

$mdFile      = "$newMdUuid.md"
$zhXlfFile   = "$newMdUuid.$newXlfHash.zh-cn.xlf"
$deXlfFile   = "$newMdUuid.$newXlfHash.de-de.xlf"

$png1Uuid    = "643b4468-45ed-49c8-b640-96203b9cdfec"
$png2Uuid    = "e67ff8ea-a00c-403c-904b-525570704630"
$png1File    = "$png1Uuid.png"
$png2File    = "$png2Uuid.png"
$png1Target  = "4c2b7c6d67b5bbb2a908e58386817f521ba28982.png"
$png2Target  = "dcdb25cbbe36e1dc51899005bb14c9cc1609c0a3.png"

$handoffDate   = "2016-03-22 13:09:40"
$targetDateZh  = "2016-03-22 13:09:36"
$epoch         = "0001-01-01 00:00:00"
$readyStatus   = "Ready for handoff"
$includeStatus = "Include"
$dependStatus  = "IsDependency"
$dependFrom    = "e2e\$mdFile"

$mdRepoUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/550209aa69855b6ed9c0a5a4403bd7414b5d7b35/e2e"
$zhXlfUrl   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/da997ef4708a58f45af33727e2e059bbf129e25d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht"
$deXlfUrl   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/38089eefdaac2350e15b093a6d3a904ea79f3ec9/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht"
$pngBaseUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/da997ef4708a58f45af33727e2e059bbf129e25d/ol-handoff/OpenLocalizationTestOrg/oltest/ci/ht"

function Set-HyperlinkFormat($range) {
    $range.Font.Underline = $true
    $range.Font.Color = 15570276   # RGB(0x64,0x95,0xED) -> BGR long
}

function Set-DateFormat($range) {
    $range.NumberFormat = "yyyy-mm-dd HH:mm:ss"
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Hyperlinks.Delete()

$ovw.Range("A2").Value = $mdFile
$ovw.Range("B2").Value = $readyStatus
$ovw.Range("C2").Value = $readyStatus
Set-DateFormat($ovw.Range("D2"))
$ovw.Range("D2").Value = $handoffDate

$ovw.Range("A3").Value = $png1File
$ovw.Range("B3").Value = $readyStatus
$ovw.Range("C3").Value = $readyStatus
Set-DateFormat($ovw.Range("D3"))
$ovw.Range("D3").Value = $handoffDate

$ovw.Range("A4").Value = $png2File
$ovw.Range("B4").Value = $readyStatus
$ovw.Range("C4").Value = $readyStatus
Set-DateFormat($ovw.Range("D4"))
$ovw.Range("D4").Value = $handoffDate

Set-HyperlinkFormat($ovw.Range("A2"))
Set-HyperlinkFormat($ovw.Range("A3"))
Set-HyperlinkFormat($ovw.Range("A4"))

$ovw.Hyperlinks.Add($ovw.Range("A2"), "$mdRepoUrl/$mdFile", "", "", $mdFile) | Out-Null
$ovw.Hyperlinks.Add($ovw.Range("A3"), "$pngBaseUrl/$png1File", "", "", $png1File) | Out-Null
$ovw.Hyperlinks.Add($ovw.Range("A4"), "$pngBaseUrl/$png2File", "", "", $png2File) | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Hyperlinks.Delete()

$zh.Range("A2").Value = $mdFile
$zh.Range("B2").Value = ".md"
$zh.Range("C2").Value = $readyStatus
$zh.Range("D2").Value = $zhXlfFile
Set-DateFormat($zh.Range("E2"))
$zh.Range("E2").Value = $targetDateZh
Set-DateFormat($zh.Range("H2"))
$zh.Range("H2").Value = $epoch
$zh.Range("J2").Value = $includeStatus

$zh.Range("A3").Value = $png1File
$zh.Range("B3").Value = ".png"
$zh.Range("C3").Value = $readyStatus
$zh.Range("D3").Value = $png1Target
Set-DateFormat($zh.Range("E3"))
$zh.Range("E3").Value = $targetDateZh
Set-DateFormat($zh.Range("H3"))
$zh.Range("H3").Value = $epoch
$zh.Range("J3").Value = $dependStatus
$zh.Range("K3").Value = $dependFrom

$zh.Range("A4").Value = $png2File
$zh.Range("B4").Value = ".png"
$zh.Range("C4").Value = $readyStatus
$zh.Range("D4").Value = $png2Target
Set-DateFormat($zh.Range("E4"))
$zh.Range("E4").Value = $targetDateZh
Set-DateFormat($zh.Range("H4"))
$zh.Range("H4").Value = $epoch
$zh.Range("J4").Value = $dependStatus
$zh.Range("K4").Value = $dependFrom

Set-HyperlinkFormat($zh.Range("A2"))
Set-HyperlinkFormat($zh.Range("D2"))
Set-HyperlinkFormat($zh.Range("A3"))
Set-HyperlinkFormat($zh.Range("D3"))
Set-HyperlinkFormat($zh.Range("A4"))
Set-HyperlinkFormat($zh.Range("D4"))

$zh.Hyperlinks.Add($zh.Range("A2"), "$mdRepoUrl/$mdFile", "", "", $mdFile) | Out-Null
$zh.Hyperlinks.Add($zh.Range("D2"), "$zhXlfUrl/$zhXlfFile", "", "", $zhXlfFile) | Out-Null
$zh.Hyperlinks.Add($zh.Range("A3"), "$pngBaseUrl/$png1File", "", "", $png1File) | Out-Null
$zh.Hyperlinks.Add($zh.Range("D3"), "$zhXlfUrl/$png1Target", "", "", $png1Target) | Out-Null
$zh.Hyperlinks.Add($zh.Range("A4"), "$pngBaseUrl/$png2File", "", "", $png2File) | Out-Null
$zh.Hyperlinks.Add($zh.Range("D4"), "$zhXlfUrl/$png2Target", "", "", $png2Target) | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Hyperlinks.Delete()

$de.Range("A2").Value = $mdFile
$de.Range("B2").Value = ".md"
$de.Range("C2").Value = $readyStatus
$de.Range("D2").Value = $deXlfFile
Set-DateFormat($de.Range("E2"))
$de.Range("E2").Value = $handoffDate
Set-DateFormat($de.Range("H2"))
$de.Range("H2").Value = $epoch
$de.Range("J2").Value = $includeStatus

$de.Range("A3").Value = $png1File
$de.Range("B3").Value = ".png"
$de.Range("C3").Value = $readyStatus
$de.Range("D3").Value = $png1Target
Set-DateFormat($de.Range("E3"))
$de.Range("E3").Value = $handoffDate
Set-DateFormat($de.Range("H3"))
$de.Range("H3").Value = $epoch
$de.Range("J3").Value = $dependStatus
$de.Range("K3").Value = $dependFrom

$de.Range("A4").Value = $png2File
$de.Range("B4").Value = ".png"
$de.Range("C4").Value = $readyStatus
$de.Range("D4").Value = $png2Target
Set-DateFormat($de.Range("E4"))
$de.Range("E4").Value = $handoffDate
Set-DateFormat($de.Range("H4"))
$de.Range("H4").Value = $epoch
$de.Range("J4").Value = $dependStatus
$de.Range("K4").Value = $dependFrom

Set-HyperlinkFormat($de.Range("A2"))
Set-HyperlinkFormat($de.Range("D2"))
Set-HyperlinkFormat($de.Range("A3"))
Set-HyperlinkFormat($de.Range("D3"))
Set-HyperlinkFormat($de.Range("A4"))
Set-HyperlinkFormat($de.Range("D4"))

$de.Hyperlinks.Add($de.Range("A2"), "$mdRepoUrl/$mdFile", "", "", $mdFile) | Out-Null
$de.Hyperlinks.Add($de.Range("D2"), "$deXlfUrl/$deXlfFile", "", "", $deXlfFile) | Out-Null
$de.Hyperlinks.Add($de.Range("A3"), "$pngBaseUrl/$png1File", "", "", $png1File) | Out-Null
$de.Hyperlinks.Add($de.Range("D3"), "$deXlfUrl/$png1Target", "", "", $png1Target) | Out-Null
$de.Hyperlinks.Add($de.Range("A4"), "$pngBaseUrl/$png2File", "", "", $png2File) | Out-Null
$de.Hyperlinks.Add($de.Range("D4"), "$deXlfUrl/$png2Target", "", "", $png2Target) | Out-Null
